$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-04-02 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-03 Thursday", 2)

$table = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

# Row 1 (table row index 1)
Set-CellText $table 1 1 "53÷2=26, 1"
Set-CellText $table 1 2 "52÷5=10, 2"
Set-CellText $table 1 3 "16÷4=4, 0"
Set-CellText $table 1 4 "86÷2=43, 0"
Set-CellText $table 1 5 "53÷3=17, 2"

# Row 5 (table row index 5, since there are 4 blank rows between data rows)
Set-CellText $table 5 1 "64÷2=32, 0"
Set-CellText $table 5 2 "19÷3=6, 1"
Set-CellText $table 5 3 "10÷3=3, 1"
Set-CellText $table 5 4 "82÷8=10, 2"
Set-CellText $table 5 5 "65÷7=9, 2"

# Row 9
Set-CellText $table 9 1 "60÷2=30, 0"
Set-CellText $table 9 2 "64÷7=9, 1"
Set-CellText $table 9 3 "26÷5=5, 1"
Set-CellText $table 9 4 "36÷8=4, 4"
Set-CellText $table 9 5 "29÷4=7, 1"

# Row 13
Set-CellText $table 13 1 "78÷6=13, 0"
Set-CellText $table 13 2 "95÷8=11, 7"
Set-CellText $table 13 3 "80÷3=26, 2"
Set-CellText $table 13 4 "59÷7=8, 3"
Set-CellText $table 13 5 "26÷4=6, 2"

# Row 17
Set-CellText $table 17 1 "93÷7=13, 2"
Set-CellText $table 17 2 "73÷3=24, 1"
Set-CellText $table 17 3 "58÷4=14, 2"
Set-CellText $table 17 4 "98÷5=19, 3"
Set-CellText $table 17 5 "75÷5=15, 0"
